# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) holds freshly recomputed strike-count values (s_vals)
# for each saved game row. Write the new values in place, row by row,
# leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,2,2,1,1,1,0,1,1,1,1,2,1,2,1,2,2,2,1,0,3,1,2,0,1,1,1,1,2,1,2,1,3,0,1,0,0,2,2,2,1,0,1,1,0,3,0,1,1,1,2,4,1,0,0,2,0,3,1,0,1,1,2,1,0,0,0,2,1,1,3,2,2,3,2,0,2)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
